$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Learning Factor and Input tier ")

# Shift existing data: row2 col A (value 4) moves from A2 to A2 (label), B2 stays as value 4
# New layout:
#  A1 = "Коэффициент скорости обучения" (label), B1 = 0.1 (existing value moved from A1)
#  A2 = "Количество входов нейронной сети" (label), B2 = 4   (existing value moved from A2)
#  A3 = "Размерность выходного слоя" (label, new row), B3 = 2 (existing value moved from B2)

$ws.Range("B1").Value = 0.1
$ws.Range("B2").Value = 4
$ws.Range("B3").Value = 2

$ws.Range("A1").Value = "Коэффициент скорости обучения"
$ws.Range("A2").Value = "Количество входов нейронной сети"
$ws.Range("A3").Value = "Размерность выходного слоя"

$ws.Range("G8").Select()
